$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Update IP addresses in column F (rows 2-6): 127.0.0.1 -> 192.168.1.113..117
$ws.Range("F2").Value = "192.168.1.113"
$ws.Range("F3").Value = "192.168.1.114"
$ws.Range("F4").Value = "192.168.1.115"
$ws.Range("F5").Value = "192.168.1.116"
$ws.Range("F6").Value = "192.168.1.117"

# Auto-fit column F to the new (longer) IP address content (best-fit width ~15 chars)
$ws.Columns.Item(6).ColumnWidth = 14.29

# Clear the contents of row 7 (the extra GameServer_2 row), keeping styles
$ws.Range("A7:H7").ClearContents()

# Update the selection to the now-empty row 7
$ws.Range("A7:XFD7").Select()
